# Updates the crypto price/volume table to the latest snapshot
# (GitHub Actions refresh, Sat Jun  3 02:53:01 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.172.38"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "1.901.95"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'306.66"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "'0.5223"
$ws.Range("E7").Value = "  +1.68%  "

$ws.Range("D8").Value = "'0.3765"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").Value = "'0.07239"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("D10").Value = "'21.16"
$ws.Range("E10").Value = "  +2.25%  "

$ws.Range("D11").Value = "'0.8978"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").Value = "'0.08430"
$ws.Range("E12").Value = "  +11.80%  "

$ws.Range("D13").Value = "1.917.32"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "'94.61"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "'5.262"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").Value = "'0.000008588"
$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").Value = "'14.49"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").Value = "27.225.50"
$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("D21").Value = "'5.055"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").Value = "2.155.35"
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  +1.71%  "

$ws.Range("D24").Value = "'6.416"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'146.66"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D26").Value = "'2.272"
$ws.Range("E26").Value = "  +7.71%  "

$ws.Range("D27").Value = "'1.752"
$ws.Range("E27").Value = "  -1.57%  "

$ws.Range("D28").Value = "'18.13"
$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").Value = "'4.919"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").Value = "'4.780"
$ws.Range("E31").Value = "  +0.57%  "

$ws.Range("D32").Value = "'0.09209"
$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("D33").Value = "'0.8102"
$ws.Range("E33").Value = "  +7.92%  "

$ws.Range("D34").Value = "'0.05058"
$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("D35").Value = "'1.237"
$ws.Range("E35").Value = "  +5.39%  "

$ws.Range("D36").Value = "'2.978"

$ws.Range("D37").Value = "'3.364"
$ws.Range("E37").Value = "  +3.49%  "

$ws.Range("D38").Value = "'2.553"
$ws.Range("E38").Value = "  +2.80%  "

$ws.Range("D39").Value = "'0.5697"
$ws.Range("E39").Value = "  +2.16%  "

$ws.Range("D40").Value = "'0.01975"
$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("D41").Value = "'1.070"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.609"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'8.940"
$ws.Range("E43").Value = "  +2.96%  "

$ws.Range("D44").Value = "'118.15"
$ws.Range("E44").Value = "  +2.39%  "

$ws.Range("D45").Value = "'0.1509"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("D46").Value = "'0.4823"
$ws.Range("E46").Value = "  +1.06%  "

$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").Value = "'10.13"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").Value = "'1.609"
$ws.Range("E49").Value = "  +2.74%  "

$ws.Range("D50").Value = "'37.38"
$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("D51").Value = "'63.55"
$ws.Range("E51").Value = "  +0.53%  "
